$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.320001602172852
$ws.Range("B1").Value = 2.744488954544067
$ws.Range("C1").Value = 3.033463954925537
$ws.Range("D1").Value = 1.521171569824219
$ws.Range("E1").Value = 1.102577090263367
